$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp label (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 26 de Abril de 2020 a las 07:52"

# --- Update Hungria row (row 65) with refreshed case counts ---
$ws.Range("B65").Value = 2500
$ws.Range("C65").Value = 57
$ws.Range("D65").Value = 485
$ws.Range("E65").Value = 1743
$ws.Range("F65").Value = 61
$ws.Range("G65").Value = 10
$ws.Range("H65").Value = 272

# --- Bulgaria moves up (now row 83) with new case counts; Ghana moves down ---
# (now row 84) keeping its previous figures unchanged.
$ws.Range("A83").Value = "Bulgaria"
$ws.Range("B83").Value = 1290
$ws.Range("C83").Value = 43
$ws.Range("D83").Value = 205
$ws.Range("E83").Value = 1030
$ws.Range("F83").Value = 37
$ws.Range("G83").Value = 0
$ws.Range("H83").Value = 55

$ws.Range("A84").Value = "Ghana"
$ws.Range("B84").Value = 1279
$ws.Range("C84").Value = 0
$ws.Range("D84").Value = 134
$ws.Range("E84").Value = 1135
$ws.Range("F84").Value = 4
$ws.Range("G84").Value = 0
$ws.Range("H84").Value = 10
